$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: materialise row 131 with the correct style by duplicating row
#     130 (still holding its original values/style at this point), then
#     overwrite every cell with the real data for the new match. ---
$ws.Range("A130:V130").Copy($ws.Range("A131"))

$ws.Cells.Item(131, 1).Value  = 130
$ws.Cells.Item(131, 2).Value  = "romania"
$ws.Cells.Item(131, 3).Value  = "liga-2"
$ws.Cells.Item(131, 4).Value  = "2023-2024"
$ws.Cells.Item(131, 5).Value  = 45243.6875
$ws.Cells.Item(131, 6).Value  = "Steaua Bucuresti"
$ws.Cells.Item(131, 7).Value  = 1
$ws.Cells.Item(131, 8).Value  = "Csikszereda M. Ciuc"
$ws.Cells.Item(131, 9).Value  = 0
$ws.Cells.Item(131, 10).Value = 1.85
$ws.Cells.Item(131, 11).Value = "09/11/2023 22:12"
$ws.Cells.Item(131, 12).Value = 1.76
$ws.Cells.Item(131, 13).Value = "13/11/2023 16:26"
$ws.Cells.Item(131, 14).Value = 3.3
$ws.Cells.Item(131, 15).Value = "09/11/2023 22:12"
$ws.Cells.Item(131, 16).Value = 3.57
$ws.Cells.Item(131, 17).Value = "13/11/2023 16:26"
$ws.Cells.Item(131, 18).Value = 3.87
$ws.Cells.Item(131, 19).Value = "09/11/2023 22:12"
$ws.Cells.Item(131, 20).Value = 4.76
$ws.Cells.Item(131, 21).Value = "13/11/2023 16:26"
$ws.Cells.Item(131, 22).Value = "https://www.betexplorer.com/football/romania/liga-2/csa-steaua-bucuresti-miercurea-ciuc/vN7d6NsT/"

# --- Step 2: rows 125 and 130 swap match data (everything but Indice/pais/
#     torneio/temporada/data_partida/home_ft_gols, i.e. columns A-E and G
#     stay put). Capture row 125's current ("Tunari vs Concordia") values
#     via .Text before overwriting, since row 130 needs them, and plain
#     .Value reads back as an rvalue do not resolve on this host. ---
$f125 = $ws.Range("F125").Text
$h125 = $ws.Range("H125").Text
$i125 = $ws.Range("I125").Text
$j125 = $ws.Range("J125").Text
$l125 = $ws.Range("L125").Text
$m125 = $ws.Range("M125").Text
$n125 = $ws.Range("N125").Text
$p125 = $ws.Range("P125").Text
$q125 = $ws.Range("Q125").Text
$r125 = $ws.Range("R125").Text
$t125 = $ws.Range("T125").Text
$u125 = $ws.Range("U125").Text
$v125 = $ws.Range("V125").Text

# Row 125 becomes the "Unirea Dej vs CSM Resita" match (row 130's old data)
$ws.Cells.Item(125, 6).Value  = "Unirea Dej"
$ws.Cells.Item(125, 8).Value  = "CSM Resita"
$ws.Cells.Item(125, 9).Value  = 3
$ws.Cells.Item(125, 10).Value = 2.53
$ws.Cells.Item(125, 12).Value = 2.51
$ws.Cells.Item(125, 13).Value = "11/11/2023 09:59"
$ws.Cells.Item(125, 14).Value = 2.97
$ws.Cells.Item(125, 16).Value = 3.12
$ws.Cells.Item(125, 17).Value = "11/11/2023 09:58"
$ws.Cells.Item(125, 18).Value = 2.7
$ws.Cells.Item(125, 20).Value = 2.96
$ws.Cells.Item(125, 21).Value = "11/11/2023 09:59"
$ws.Cells.Item(125, 22).Value = "https://www.betexplorer.com/football/romania/liga-2/unirea-dej-csm-resita/j1yFO4cc/"

# Row 130 becomes the "Tunari vs Concordia" match (row 125's old data)
$ws.Cells.Item(130, 6).Value  = $f125
$ws.Cells.Item(130, 8).Value  = $h125
$ws.Cells.Item(130, 9).Value  = $i125
$ws.Cells.Item(130, 10).Value = $j125
$ws.Cells.Item(130, 12).Value = $l125
$ws.Cells.Item(130, 13).Value = $m125
$ws.Cells.Item(130, 14).Value = $n125
$ws.Cells.Item(130, 16).Value = $p125
$ws.Cells.Item(130, 17).Value = $q125
$ws.Cells.Item(130, 18).Value = $r125
$ws.Cells.Item(130, 20).Value = $t125
$ws.Cells.Item(130, 21).Value = $u125
$ws.Cells.Item(130, 22).Value = $v125
